$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $origStyle = $Cell.Style
    $Cell.Formula = "'" + $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '319.07'  # D2: 319.05 -> 319.07
Set-TextValue $ws.Cells.Item(2, 5) '3.49%'  # E2: 3.30% -> 3.49%
Set-TextValue $ws.Cells.Item(3, 4) '41.47'  # D3: 41.45 -> 41.47
Set-TextValue $ws.Cells.Item(3, 5) '1.41%'  # E3: 1.23% -> 1.41%
Set-TextValue $ws.Cells.Item(4, 5) '2.28%'  # E4: 2.33% -> 2.28%
Set-TextValue $ws.Cells.Item(5, 5) '1.49%'  # E5: 1.52% -> 1.49%
Set-TextValue $ws.Cells.Item(6, 2) 'FTXToken'  # B6: GateToken -> FTXToken
Set-TextValue $ws.Cells.Item(6, 3) 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'  # C6: https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt -> https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt
Set-TextValue $ws.Cells.Item(6, 4) '1.698'  # D6: 4.321 -> 1.698
Set-TextValue $ws.Cells.Item(6, 5) '4.65%'  # E6: 1.01% -> 4.65%
Set-TextValue $ws.Cells.Item(7, 2) 'MXToken'  # B7: FTXToken -> MXToken
Set-TextValue $ws.Cells.Item(7, 3) 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'  # C7: https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt -> https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx
Set-TextValue $ws.Cells.Item(7, 4) '0.9504'  # D7: 1.698 -> 0.9504
Set-TextValue $ws.Cells.Item(7, 5) '4.50%'  # E7: 4.90% -> 4.50%
Set-TextValue $ws.Cells.Item(8, 2) 'BTSEToken'  # B8: MXToken -> BTSEToken
Set-TextValue $ws.Cells.Item(8, 3) 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'  # C8: https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx -> https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse
Set-TextValue $ws.Cells.Item(8, 4) '2.425'  # D8: 0.9521 -> 2.425
Set-TextValue $ws.Cells.Item(8, 5) '-2.52%'  # E8: 4.60% -> -2.52%
Set-TextValue $ws.Cells.Item(9, 2) 'LiechtensteinCryptoassetsExchange'  # B9: BTSEToken -> LiechtensteinCryptoassetsExchange
Set-TextValue $ws.Cells.Item(9, 3) 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'  # C9: https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse -> https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx
Set-TextValue $ws.Cells.Item(9, 4) '0.1245'  # D9: 2.425 -> 0.1245
Set-TextValue $ws.Cells.Item(9, 5) '5.35%'  # E9: -2.52% -> 5.35%
Set-TextValue $ws.Cells.Item(10, 2) 'WazirX'  # B10: LiechtensteinCryptoassetsExchange -> WazirX
Set-TextValue $ws.Cells.Item(10, 3) 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'  # C10: https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx -> https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx
Set-TextValue $ws.Cells.Item(10, 4) '0.1833'  # D10: 0.1261 -> 0.1833
Set-TextValue $ws.Cells.Item(10, 5) '1.03%'  # E10: 6.42% -> 1.03%
Set-TextValue $ws.Cells.Item(11, 2) 'MandalaExchangeToken'  # B11: WazirX -> MandalaExchangeToken
Set-TextValue $ws.Cells.Item(11, 3) 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'  # C11: https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx -> https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx
Set-TextValue $ws.Cells.Item(11, 4) '0.09195'  # D11: 0.1834 -> 0.09195
Set-TextValue $ws.Cells.Item(11, 5) '1.98%'  # E11: 0.73% -> 1.98%
Set-TextValue $ws.Cells.Item(12, 2) 'BitrueCoin'  # B12: MandalaExchangeToken -> BitrueCoin
Set-TextValue $ws.Cells.Item(12, 3) 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'  # C12: https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx -> https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr
Set-TextValue $ws.Cells.Item(12, 4) '0.04367'  # D12: 0.09168 -> 0.04367
Set-TextValue $ws.Cells.Item(12, 5) '2.15%'  # E12: 0.60% -> 2.15%
Set-TextValue $ws.Cells.Item(13, 2) 'BitMartToken'  # B13: BitrueCoin -> BitMartToken
Set-TextValue $ws.Cells.Item(13, 3) 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'  # C13: https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr -> https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx
Set-TextValue $ws.Cells.Item(13, 4) '0.1051'  # D13: 0.04370 -> 0.1051
Set-TextValue $ws.Cells.Item(13, 5) '0.48%'  # E13: 2.05% -> 0.48%
Set-TextValue $ws.Cells.Item(14, 2) 'BitForexToken'  # B14: BitMartToken -> BitForexToken
Set-TextValue $ws.Cells.Item(14, 3) 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'  # C14: https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx -> https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf
Set-TextValue $ws.Cells.Item(14, 4) '0.001283'  # D14: 0.1050 -> 0.001283
Set-TextValue $ws.Cells.Item(14, 5) '2.40%'  # E14: 0.50% -> 2.40%
Set-TextValue $ws.Cells.Item(15, 2) 'TigerCash'  # B15: BitForexToken -> TigerCash
Set-TextValue $ws.Cells.Item(15, 3) 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'  # C15: https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf -> https://coinranking.com/coin/6hIn06L2+tigercash-tch
Set-TextValue $ws.Cells.Item(15, 4) '0.005977'  # D15: 0.001288 -> 0.005977
Set-TextValue $ws.Cells.Item(15, 5) '2.05%'  # E15: 1.98% -> 2.05%
Set-TextValue $ws.Cells.Item(16, 2) 'LEO'  # B16: TigerCash -> LEO
Set-TextValue $ws.Cells.Item(16, 3) 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'  # C16: https://coinranking.com/coin/6hIn06L2+tigercash-tch -> https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo
Set-TextValue $ws.Cells.Item(16, 4) '3.340'  # D16: 0.005987 -> 3.340
Set-TextValue $ws.Cells.Item(16, 5) '-0.34%'  # E16: 1.71% -> -0.34%
Set-TextValue $ws.Cells.Item(17, 2) 'GateToken'  # B17: LEO -> GateToken
Set-TextValue $ws.Cells.Item(17, 3) 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'  # C17: https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo -> https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt
Set-TextValue $ws.Cells.Item(17, 4) '4.321'  # D17: 3.337 -> 4.321
Set-TextValue $ws.Cells.Item(17, 5) '1.01%'  # E17: -0.41% -> 1.01%
Set-TextValue $ws.Cells.Item(18, 4) '0.3342'  # D18: 0.3353 -> 0.3342
Set-TextValue $ws.Cells.Item(18, 5) '2.60%'  # E18: 2.91% -> 2.60%
Set-TextValue $ws.Cells.Item(19, 4) '7.690'  # D19: 7.657 -> 7.690
Set-TextValue $ws.Cells.Item(19, 5) '11.24%'  # E19: 10.91% -> 11.24%
Set-TextValue $ws.Cells.Item(20, 5) '-3.69%'  # E20: -4.19% -> -3.69%
Set-TextValue $ws.Cells.Item(22, 5) '-0.68%'  # E22: -0.62% -> -0.68%
Set-TextValue $ws.Cells.Item(23, 4) '0.001264'  # D23: 0.001263 -> 0.001264
Set-TextValue $ws.Cells.Item(23, 5) '-0.72%'  # E23: -0.64% -> -0.72%
Set-TextValue $ws.Cells.Item(24, 4) '0.004121'  # D24: 0.004124 -> 0.004121
Set-TextValue $ws.Cells.Item(24, 5) '-0.09%'  # E24: 0.06% -> -0.09%
Set-TextValue $ws.Cells.Item(25, 5) '-0.20%'  # E25: -0.35% -> -0.20%
Set-TextValue $ws.Cells.Item(38, 4) '0.02543'  # D38: 0.02557 -> 0.02543
Set-TextValue $ws.Cells.Item(38, 5) '5.43%'  # E38: 5.60% -> 5.43%
Set-TextValue $ws.Cells.Item(39, 4) '0.05346'  # D39: 0.05359 -> 0.05346
Set-TextValue $ws.Cells.Item(39, 5) '2.48%'  # E39: 2.43% -> 2.48%
Set-TextValue $ws.Cells.Item(40, 4) '0.007777'  # D40: 0.007755 -> 0.007777
Set-TextValue $ws.Cells.Item(40, 5) '-0.63%'  # E40: -0.40% -> -0.63%
Set-TextValue $ws.Cells.Item(41, 4) '0.1319'  # D41: 0.1318 -> 0.1319
Set-TextValue $ws.Cells.Item(41, 5) '1.49%'  # E41: 1.46% -> 1.49%
Set-TextValue $ws.Cells.Item(42, 5) '7.75%'  # E42: 7.76% -> 7.75%
Set-TextValue $ws.Cells.Item(43, 5) '2.86%'  # E43: 2.75% -> 2.86%
Set-TextValue $ws.Cells.Item(44, 4) '0.007582'  # D44: 0.007571 -> 0.007582
Set-TextValue $ws.Cells.Item(44, 5) '-6.16%'  # E44: -6.17% -> -6.16%
Set-TextValue $ws.Cells.Item(45, 4) '0.3443'  # D45: 0.3441 -> 0.3443
Set-TextValue $ws.Cells.Item(45, 5) '12.44%'  # E45: 11.96% -> 12.44%
Set-TextValue $ws.Cells.Item(46, 4) '0.00006686'  # D46: 0.00006684 -> 0.00006686
Set-TextValue $ws.Cells.Item(46, 5) '-3.15%'  # E46: -3.09% -> -3.15%
Set-TextValue $ws.Cells.Item(47, 5) '-0.25%'  # E47: -0.24% -> -0.25%
Set-TextValue $ws.Cells.Item(48, 4) '0.2184'  # D48: 0.1806 -> 0.2184
Set-TextValue $ws.Cells.Item(48, 5) '127.69%'  # E48: 83.34% -> 127.69%
Set-TextValue $ws.Cells.Item(49, 5) '39.77%'  # E49: 39.76% -> 39.77%
Set-TextValue $ws.Cells.Item(50, 5) '-0.25%'  # E50: -0.24% -> -0.25%
Set-TextValue $ws.Cells.Item(51, 5) '-0.25%'  # E51: -0.24% -> -0.25%
